$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.192.50"
$ws.Range("E2").Value = "  -4.67%  "

# Row 3
$ws.Range("D3").Value = "2.540.52"
$ws.Range("E3").Value = "  -4.12%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'504.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.30%  "

# Row 6
$ws.Range("D6").Value = "'144.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.60%  "

# Row 7
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "'0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.54%  "

# Row 9
$ws.Range("D9").Value = "2.544.94"
$ws.Range("E9").Value = "  -4.47%  "

# Row 10
$ws.Range("D10").Value = "'6.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.21%  "

# Row 11
$ws.Range("D11").Value = "'0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.07%  "

# Row 12
$ws.Range("D12").Value = "'0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.02%  "

# Row 13
$ws.Range("E13").Value = "  -0.76%  "

# Row 14
$ws.Range("D14").Value = "2.976.28"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15
$ws.Range("D15").Value = "58.146.37"
$ws.Range("E15").Value = "  -4.75%  "

# Row 16
$ws.Range("D16").Value = "'20.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.33%  "

# Row 17
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.28%  "

# Row 18
$ws.Range("D18").Value = "2.536.48"
$ws.Range("E18").Value = "  -4.05%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'343.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.37%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.82%  "

# Row 21
$ws.Range("D21").Value = "'10.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.49%  "

# Row 22
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").Value = "'5.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.93%  "

# Row 24
$ws.Range("D24").Value = "'60.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "

# Row 25
$ws.Range("D25").Value = "'0.409"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.16%  "

# Row 26
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "

# Row 27
$ws.Range("D27").Value = "2.647.82"
$ws.Range("E27").Value = "  -4.40%  "

# Row 28
$ws.Range("D28").Value = "'0.158"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.36%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0785"
$ws.Range("E29").Value = "  -8.58%  "

# Row 30
$ws.Range("D30").Value = "'6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.32%  "

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("D32").Value = "'5.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.31%  "

# Row 33
$ws.Range("D33").Value = "'149.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "

# Row 34
$ws.Range("D34").Value = "'18.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.66%  "

# Row 35
$ws.Range("D35").Value = "'1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.69%  "

# Row 36
$ws.Range("D36").Value = "'3.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.95%  "

# Row 37
$ws.Range("D37").Value = "'0.903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "

# Row 38
$ws.Range("D38").Value = "'1.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.74%  "

# Row 39
$ws.Range("D39").Value = "'35.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "

# Row 40
$ws.Range("D40").Value = "'0.824"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.60%  "

# Row 41
$ws.Range("D41").Value = "'1.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.23%  "

# Row 42
$ws.Range("D42").Value = "'3.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.00%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'281.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.88%  "

# Row 44
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.17%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("D46").Value = "'0.596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.88%  "

# Row 47
$ws.Range("D47").Value = "'0.0531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.81%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'18.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.13%  "

# Row 49
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50
$ws.Range("D50").Value = "'0.0226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.63%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.99%  "

